$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the shared string "R40". Retype it as the text "1"
# (a leading apostrophe forces Excel to store it as text rather than the
# number 1, so the cell keeps its shared-string ("t=s") representation
# instead of becoming numeric).
$ws.Range("B11").Formula = "'1"
